$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 833.0395
$ws.Range("I15").Value = 833.0395
$ws.Range("K15").Value = 2499.1185
$ws.Range("M15").Value = -2330.1185

$ws.Range("H16").Value = 20750
$ws.Range("I16").Value = 25000
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -24770
$ws.Range("N16").Value = -8460

$ws.Range("H64").Value = 19668.23
$ws.Range("I64").Value = 20807.25
$ws.Range("K64").Value = 20807.25
$ws.Range("M64").Value = -20559.25

$ws.Range("H67").Value = 19668.23
$ws.Range("I67").Value = 20807.25
$ws.Range("K67").Value = 20807.25
$ws.Range("M67").Value = -19949.25

$ws.Range("H70").Value = 3591278.8
$ws.Range("I70").Value = 8548971
$ws.Range("K70").Value = 25646913
$ws.Range("M70").Value = -25646643

$ws.Range("H73").Value = 3591278.8
$ws.Range("I73").Value = 8548971
$ws.Range("K73").Value = 25646913
$ws.Range("M73").Value = -25645977

$ws.Range("H76").Value = 5498.385
$ws.Range("J76").Value = 5990.2
$ws.Range("L76").Value = 5990.2
$ws.Range("N76").Value = -6620.2

$ws.Range("H79").Value = 5498.385
$ws.Range("J79").Value = 5990.2
$ws.Range("L79").Value = 5990.2
$ws.Range("N79").Value = -8174.2

$ws.Range("H112").Value = 6863.7393
$ws.Range("J112").Value = 7161.0684
$ws.Range("L112").Value = 21483.2052
$ws.Range("N112").Value = -23699.2052

$ws.Range("H125").Value = 35800.89
$ws.Range("I125").Value = 60866.2
$ws.Range("J125").Value = 4469.25
$ws.Range("K125").Value = 547795.7999999999
$ws.Range("L125").Value = 40223.25
$ws.Range("M125").Value = -545335.7999999999
$ws.Range("N125").Value = -45143.25

$ws.Range("H127").Value = 1260.6154
$ws.Range("J127").Value = 1802.8572
$ws.Range("L127").Value = 5408.571599999999
$ws.Range("N127").Value = -15328.5716

$ws.Range("H137").Value = 10059.613
$ws.Range("I137").Value = 23450.066
$ws.Range("J137").Value = 3133.5173
$ws.Range("K137").Value = 70350.198
$ws.Range("L137").Value = 9400.5519
$ws.Range("M137").Value = -67800.198
$ws.Range("N137").Value = -14500.5519

$ws.Range("H138").Value = 2994.4595
$ws.Range("I138").Value = 2461.1738
$ws.Range("K138").Value = 7383.5214
$ws.Range("M138").Value = -2243.5214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3980.91
$ws.Range("I32").Value = 3532.883
$ws.Range("J32").Value = 11000
$ws.Range("K32").Value = 3532.883
$ws.Range("L32").Value = 11000
$ws.Range("M32").Value = -3245.883
$ws.Range("N32").Value = -11574

$ws.Range("H45").Value = 4064.4
$ws.Range("I45").Value = 2268.5454
$ws.Range("J45").Value = 5475.4287
$ws.Range("K45").Value = 2268.5454
$ws.Range("L45").Value = 5475.4287
$ws.Range("M45").Value = -1891.5454
$ws.Range("N45").Value = -6229.4287

$ws.Range("H102").Value = 4306.2856
$ws.Range("I102").Value = 3399
$ws.Range("J102").Value = 5516
$ws.Range("K102").Value = 3399
$ws.Range("L102").Value = 5516
$ws.Range("M102").Value = -1777
$ws.Range("N102").Value = -8760

$ws.Range("H110").Value = 890.73334
$ws.Range("I110").Value = 889.4167
$ws.Range("K110").Value = 889.4167
$ws.Range("M110").Value = 1155.5833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1604.5
$ws.Range("I86").Value = 1604.5
$ws.Range("K86").Value = 1604.5
$ws.Range("M86").Value = -481.5

$ws.Range("H89").Value = 1604.5
$ws.Range("I89").Value = 1604.5
$ws.Range("K89").Value = 8022.5
$ws.Range("M89").Value = -2406.5

$ws.Range("H107").Value = 22577.25
$ws.Range("I107").Value = 22577.25
$ws.Range("K107").Value = 22577.25
$ws.Range("M107").Value = -20657.25

$ws.Range("H134").Value = 2266
$ws.Range("I134").Value = 2206.3914
$ws.Range("J134").Value = 2723
$ws.Range("K134").Value = 6619.174199999999
$ws.Range("L134").Value = 8169
$ws.Range("M134").Value = -4084.174199999999
$ws.Range("N134").Value = -13239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4551536
$ws.Range("I31").Value = 7696343
$ws.Range("J31").Value = 9037
$ws.Range("K31").Value = 7696343
$ws.Range("L31").Value = 9037
$ws.Range("M31").Value = -7696048
$ws.Range("N31").Value = -9627

$ws.Range("H34").Value = 4551536
$ws.Range("I34").Value = 7696343
$ws.Range("J34").Value = 9037
$ws.Range("K34").Value = 7696343
$ws.Range("L34").Value = 9037
$ws.Range("M34").Value = -7696141
$ws.Range("N34").Value = -9441

$ws.Range("H58").Value = 3053.2942
$ws.Range("J58").Value = 1309.6
$ws.Range("L58").Value = 1309.6
$ws.Range("N58").Value = -1715.6

$ws.Range("H107").Value = 475.81818
$ws.Range("I107").Value = 269.30768
$ws.Range("J107").Value = 774.1111
$ws.Range("K107").Value = 269.30768
$ws.Range("L107").Value = 774.1111
$ws.Range("M107").Value = 1650.69232
$ws.Range("N107").Value = -4614.1111

$ws.Range("H134").Value = 2137
$ws.Range("I134").Value = 1646.3334
$ws.Range("K134").Value = 4939.0002
$ws.Range("M134").Value = -2404.0002

$ws.Range("H136").Value = 3053.2942
$ws.Range("J136").Value = 1309.6
$ws.Range("L136").Value = 3928.8
$ws.Range("N136").Value = -9028.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1902.4839
$ws.Range("J113").Value = 1735.4
$ws.Range("L113").Value = 5206.200000000001
$ws.Range("N113").Value = -9546.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8015
$ws.Range("J80").Value = 18633.334
$ws.Range("L80").Value = 18633.334
$ws.Range("N80").Value = -20629.334

$ws.Range("H83").Value = 8015
$ws.Range("J83").Value = 18633.334
$ws.Range("L83").Value = 93166.67
$ws.Range("N83").Value = -103150.67

$ws.Range("H122").Value = 3466.524
$ws.Range("I122").Value = 3173.7896
$ws.Range("K122").Value = 9521.3688
$ws.Range("M122").Value = -7071.3688

$ws.Range("H132").Value = 1527.8823
$ws.Range("I132").Value = 1516
$ws.Range("K132").Value = 4548
$ws.Range("M132").Value = -2018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7306.1665
$ws.Range("I46").Value = 950.5
$ws.Range("K46").Value = 950.5
$ws.Range("M46").Value = -762.5

$ws.Range("H68").Value = 5681.25
$ws.Range("I68").Value = 2975
$ws.Range("K68").Value = 2975
$ws.Range("M68").Value = -2226

$ws.Range("H71").Value = 5681.25
$ws.Range("I71").Value = 2975
$ws.Range("K71").Value = 14875
$ws.Range("M71").Value = -11131

$ws.Range("H122").Value = 18844.154
$ws.Range("J122").Value = 19998.5
$ws.Range("L122").Value = 59995.5
$ws.Range("N122").Value = -64895.5

$ws.Range("H136").Value = 3114
$ws.Range("I136").Value = 2925.6924
$ws.Range("K136").Value = 8777.0772
$ws.Range("M136").Value = -6227.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 46014
$ws.Range("J40").Value = 46014
$ws.Range("L40").Value = 46014
$ws.Range("N40").Value = -46312

$ws.Range("H126").Value = 314285.3
$ws.Range("I126").Value = 1997.4615
$ws.Range("K126").Value = 5992.3845
$ws.Range("M126").Value = -3522.3845
